# "Fruta / hortaliza, semanal" - weekly roll of the Betarraga price series.
# A new week's observation is inserted at the top of the data block (row 129),
# every existing row shifts down by one, and the last row (old 238) lands in
# a brand-new row 239 at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 129..238 down into 130..239 by copying each source row (A:R)
# into the row immediately below it, walking bottom-up so we never clobber
# data before it has been copied onward.
for ($r = 238; $r -ge 129; $r--) {
    $src = $ws.Range("A" + $r + ":R" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":R" + ($r + 1))
    $src.Copy($dst)
}

# Row 129 keeps its place but now reports the newest week: new date + new
# volume, while the price columns (K/L/M/P) are unchanged.
$ws.Range("D129").Value = 44566
$ws.Range("J129").Value = 4000
